$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "SSSOM" (3rd sheet) - main data updates
# ------------------------------------------------------------------
$sssom = $wb.Worksheets.Item(3)

# Header row: "Reviewer1"/"Reviewer2" columns are now both labelled "reviewer_id"
$sssom.Range("L1").Value = "reviewer_id"
$sssom.Range("M1").Value = "reviewer_id"

# Switch the ORCID URLs to CURIE form, and the confidence value to lower-case
$sssom.Range("G2:G5").Value = "orcid:0009-0001-6090-9959"
$sssom.Range("J2:J5").Value = "high"
$sssom.Range("L2:L5").Value = "orcid:0000-0002-2568-5945"
$sssom.Range("M2:M5").Value = "orcid:0000-0003-4254-8683"

# Remove the hyperlink that used to decorate the author_id column
$sssom.Hyperlinks.Delete()

# Reset the (former) hyperlink formatting on the reviewer columns back to Normal
$sssom.Range("L2:M5").Style = "Normal"

# Extend the table with two extra (mostly empty) rows before restyling, so the
# whole author_id column (including the new blank cells) picks up one single
# replacement font/style combo instead of several incremental ones.
$sssom.Range("G6").Value = ""
$sssom.Range("G7").Value = ""
$sssom.Range("G6:G7").Style = $sssom.Range("G2").Style

# Restyle the author_id column with plain (non-hyperlink) formatting
$sssom.Range("G2:G7").Font.Name = "Arial"
$sssom.Range("G2:G7").Font.Size = 9
$sssom.Range("G2:G7").Font.Color = 0
$sssom.Range("G2:G7").Font.Underline = $false

# Remove the now-unused Hyperlink cell style
$wb.Styles.Item("Hyperlink").Delete()

# ------------------------------------------------------------------
# Sheet "OceanAccountsInfo" (1st sheet) - Source column style tweak
# ------------------------------------------------------------------
$info = $wb.Worksheets.Item(1)
$info.Range("D2:D5").Style = "Normal"
$info.Range("D2:D5").VerticalAlignment = -4160  ## xlTop

# ------------------------------------------------------------------
# Selection / active-sheet bookkeeping
# ------------------------------------------------------------------
$sssom.Activate()
$sssom.Range("J2:J5").Select()

Write-Output "done"
